# SCD0172 - Melakukan Proses Pemantauan pada Menu Pipeline
# Apply the content changes described by the diff:
#  - C2: replace the old "CIF" scenario text with the new short summary string
#  - D2: replace the old "CIF" scenario text with the "Jenis Nasabah" scenario text
#        (the text that used to live in C2)
#  - O2: clear the stray numeric value ("10186699862")
#  - sheet selection/top-left cell moves to around column E / cell O2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing "Jenis Nasabah" scenario text currently stored in C2
# before it gets overwritten, so it can be moved into D2.
$jenisNasabahText = $ws.Range("C2").Value2

# D2 previously held the "CIF" variant of the scenario text; it now takes
# the "Jenis Nasabah" variant that used to be in C2.
$ws.Range("D2").Value = $jenisNasabahText

# C2 gets a brand-new short description string.
$ws.Range("C2").Value = "Sales mengakses menu: Report Menu - Product Holding Ratio - Report"

# O2 no longer carries the stray "10186699862" value.
$ws.Range("O2").ClearContents()

# Update the view: scroll so column E is the left-most visible column and
# the active selection moves to O2.
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("O2").Select()
